# Apply the "Proprietari" sheet rework:
#  - remove the two unused, empty columns (E:F) so the "Loc Veci"
#    relational-key column shifts left from G to E
#  - make "Proprietari" the active sheet/tab (was "Operatii")
#  - leave the selection on the cell the user ended up on (F15)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Proprietari")

# Columns E and F are blank spacer columns between the "Localitate" column
# (D) and the "Loc Veci" key column (G). Deleting them shifts G -> E.
$ws.Range("E1:F1").EntireColumn.Delete()

# Switch the active tab from "Operatii" to "Proprietari".
$ws.Activate()

# Park the selection/active cell where the user left it after the edit.
$ws.Range("F15").Select()
